$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5633, 46005.95833333334),
    @(5628, 46005.96875),
    @(5572, 46005.97916666666),
    @(5529, 46005.98958333334),
    @(5518, 46006),
    @(5451, 46006.01041666666),
    @(5382, 46006.02083333334),
    @(5370, 46006.03125),
    @(5346, 46006.04166666666),
    @(5326, 46006.05208333334),
    @(5306, 46006.0625),
    @(5345, 46006.07291666666),
    @(5335, 46006.08333333334),
    @(5296, 46006.09375),
    @(5253, 46006.10416666666),
    @(5284, 46006.11458333334),
    @(5296, 46006.125),
    @(5318, 46006.13541666666),
    @(5363, 46006.14583333334),
    @(5389, 46006.15625),
    @(5489, 46006.16666666666),
    @(5589, 46006.17708333334),
    @(5717, 46006.1875),
    @(5817, 46006.19791666666),
    @(6075, 46006.20833333334),
    @(6336, 46006.21875),
    @(6537, 46006.22916666666),
    @(6791, 46006.23958333334),
    @(7161, 46006.25),
    @(7394, 46006.26041666666),
    @(7638, 46006.27083333334),
    @(7748, 46006.28125)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove now-unused rows 34-46 (previously had data, new dataset ends at row 33)
$lastRow = $data.Length + 2 - 1
$oldLastRow = 46
if ($oldLastRow -gt $lastRow) {
    $deleteRange = $ws.Range("A$($lastRow + 1):B$oldLastRow")
    $deleteRange.EntireRow.Delete()
}
